$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Helper: force a run split at a given character position by toggling Bold
# on/off across [pos, endPos). Word (and this host) re-normalises runs with
# identical rPr after any formatting/text change, so flipping Bold true then
# false is a no-op visually/semantically but makes the engine re-split the
# run boundaries at `pos` and `endPos`.
# ---------------------------------------------------------------------------
function Split-RunBoundary($pos, $endPos) {
    $rng = $d.Range($pos, $endPos)
    $rng.Font.Bold = $true
    $rng.Font.Bold = $false
}

# ===========================================================================
# Edit 1: "81" -> "79", and split the following run into " " + "genes..."
# ===========================================================================
$rFind = $d.Content
$rFind.Find.Execute("Somatic variant analysis of ") | Out-Null
$p1End = $rFind.End            # start of "81"
$rFind2 = $d.Content
$rFind2.Find.Execute("81 genes with clinical significance in haematological malignancy. Refer to Panel Summary for gene list.") | Out-Null
$paraEnd = $rFind2.End         # end of the paragraph's visible text

$numStart = $p1End
$numEnd = $numStart + 2        # "81" / "79" are both 2 characters
$spaceEnd = $numEnd + 1        # the single space after the number

# Replace the number text (this may merge adjacent same-format runs).
$d.Range($numStart, $numEnd).Text = "79"

# Re-establish the original run boundaries: before the number, after the
# number (before the space), and after the space (before "genes").
Split-RunBoundary $numStart $paraEnd
Split-RunBoundary $numEnd $paraEnd
Split-RunBoundary $spaceEnd $paraEnd

# ===========================================================================
# Edit 2: ".Gly646Trpfs*12 (detection limit ~ 5%-10%), CEBPA and TERT
#          (detection limit ~ 10%)"
#      -> ".Gly646Trpfs*12 (detection limit ~ 5%-10%)" + " and CEBPA" +
#         " (detection limit ~ 10%)"   (three separate runs)
# ===========================================================================
$rFind3 = $d.Content
$rFind3.Find.Execute(".Gly646Trpfs*12 (detection limit ~ 5%-10%)") | Out-Null
$keepEnd = $rFind3.End         # end of the unchanged leading fragment

$rFind4 = $d.Content
$rFind4.Find.Execute(", CEBPA and TERT") | Out-Null
$midStart = $rFind4.Start
$midOldEnd = $rFind4.End

$rFind5 = $d.Content
$rFind5.Find.Execute(". This assay is primarily qualitative") | Out-Null
$tailOldEnd = $rFind5.Start     # end of the whole original run (before the
                                 # ". This assay..." run that follows it)

# Perform the text substitution: ", CEBPA and TERT" -> " and CEBPA"
$d.Range($midStart, $midOldEnd).Text = " and CEBPA"

$midEnd = $midStart + (" and CEBPA").Length
$tailEnd = $tailOldEnd - ($midOldEnd - $midStart) + (" and CEBPA").Length

# Re-establish the leading/middle run boundaries (the trailing boundary,
# before the following ". This assay..." run, is unaffected by our edit and
# needs no explicit split).
Split-RunBoundary $keepEnd $tailEnd
Split-RunBoundary $midEnd $tailEnd

# ===========================================================================
# Edit 3: "16-Nov-2023" -> "17-Nov-2023" (single run, no split required)
# ===========================================================================
$rDate = $d.Content
$rDate.Find.Execute("16-Nov-2023") | Out-Null
$d.Range($rDate.Start, $rDate.End).Text = "17-Nov-2023"
